$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -11.738
$ws.Range("B7").Value = 5.515
$ws.Range("A8").Value = -22.084
$ws.Range("A10").Value = -21.223
$ws.Range("D10").Value = -8.068999999999999
$ws.Range("A12").Value = -21.593
$ws.Range("D12").Value = -7.074
$ws.Range("D13").Value = -8.554
$ws.Range("D14").Value = -7.762
$ws.Range("B15").Value = 5.243
$ws.Range("A18").Value = -21.696
$ws.Range("B18").Value = 6.422
$ws.Range("C18").Value = -12.113
$ws.Range("C19").Value = -12.051
$ws.Range("B20").Value = 6.956999999999999
$ws.Range("C27").Value = -13.103
$ws.Range("B29").Value = 4.893
$ws.Range("D29").Value = -7.170999999999999
$ws.Range("B30").Value = 5.675
$ws.Range("B31").Value = 5.049
$ws.Range("C31").Value = -13.286
$ws.Range("D32").Value = -8.274000000000001
$ws.Range("D35").Value = -7.387
$ws.Range("A37").Value = -20.02
$ws.Range("C38").Value = -12.845
$ws.Range("B40").Value = 8.901
$ws.Range("C42").Value = -12.328
$ws.Range("D43").Value = -8.248999999999999
$ws.Range("C44").Value = -12.777
$ws.Range("C47").Value = -12.085
$ws.Range("D48").Value = -7.489
$ws.Range("D49").Value = -7.997000000000002
$ws.Range("B50").Value = 5.867
$ws.Range("D50").Value = -8.006
$ws.Range("A55").Value = -21.868
$ws.Range("D56").Value = -8.134
$ws.Range("C58").Value = -12.727
$ws.Range("C65").Value = -12.282
$ws.Range("A68").Value = -21.681
$ws.Range("B68").Value = 5.348000000000001
$ws.Range("D69").Value = -7.540999999999999
$ws.Range("C73").Value = -12.77
$ws.Range("B76").Value = 6.689
$ws.Range("A77").Value = -20.252
$ws.Range("A78").Value = -19.774
$ws.Range("A81").Value = -21.785
$ws.Range("D81").Value = -7.650999999999999
$ws.Range("A82").Value = -21.767
$ws.Range("B87").Value = 5.508999999999999
$ws.Range("B88").Value = 5.816000000000001
$ws.Range("C90").Value = -13.427
$ws.Range("D92").Value = -6.636
$ws.Range("C94").Value = -11.085
$ws.Range("C95").Value = -11.928
$ws.Range("B96").Value = 7.187
$ws.Range("B98").Value = 5.599
$ws.Range("B101").Value = 7.717000000000001
$ws.Range("C101").Value = -12.978
$ws.Range("B102").Value = 8.009
